$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.456.30"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "2.589.22"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.12%  "
$ws.Range("D9").Value = "2.600.19"
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("E12").Value = "  +10.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.358"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.24%  "
$ws.Range("D14").Value = "3.045.78"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.57%  "
$ws.Range("D16").Value = "59.390.69"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "2.591.17"
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.468"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "0.0₃0780"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.885"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.875"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "134.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0977"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").Value = "1.955.73"
$ws.Range("E51").Value = "  -0.68%  "
